# Actualizacion automatica del mapa (2025-10-20 14:47:22)
# Adds new PEBCOM claim rows 106-111 to the "PEBCOM" worksheet,
# extending the sheet's used range from A1:R105 to A1:R111.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as genuine Text (matching the source file's
# inlineStr cells, e.g. numeric-looking case/OT codes, dates stored as
# plain strings, etc.) without leaving a residual number-format style
# on the cell once done - mirrors the existing un-styled data rows.
function Set-TextCell($addr, $val) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

    # Row 106
    Set-TextCell "A106" "7538"
    Set-TextCell "B106" "10/20/2025"
    Set-TextCell "C106" "MONTES DE OCA, MANUEL AV. 1795"
    Set-TextCell "D106" "4"
    Set-TextCell "E106" "810398922"
    Set-TextCell "F106" "PEBCOM"
    Set-TextCell "G106" "Pendiente"
    Set-TextCell "H106" "Picada"
    $ws.Range("I106").Value = 1
    Set-TextCell "J106" "Cambio"
    Set-TextCell "K106" "Sin equipos"
    Set-TextCell "L106" "Pasante"
    $ws.Range("M106").Value = -58.372942
    $ws.Range("N106").Value = -34.648042
    Set-TextCell "O106" "San Telmo"
    Set-TextCell "P106" "Capital Sur"
    Set-TextCell "Q106" "CON-H"
    Set-TextCell "R106" "Fuera de Poligono OVL"

    # Row 107
    Set-TextCell "A107" "7541"
    Set-TextCell "B107" "10/20/2025"
    Set-TextCell "C107" "MONTES DE OCA, MANUEL AV. 1245"
    Set-TextCell "D107" "4"
    Set-TextCell "E107" "810398924"
    Set-TextCell "F107" "PEBCOM"
    Set-TextCell "G107" "Pendiente"
    Set-TextCell "H107" "Picada"
    $ws.Range("I107").Value = 1
    Set-TextCell "J107" "Cambio"
    Set-TextCell "K107" "Sin equipos"
    Set-TextCell "L107" "Pasante"
    $ws.Range("M107").Value = -58.373951
    $ws.Range("N107").Value = -34.642736
    Set-TextCell "O107" "San Telmo"
    Set-TextCell "P107" "Capital Sur"
    Set-TextCell "Q107" "CON-A"
    Set-TextCell "R107" "Fuera de Poligono OVL"

    # Row 108
    Set-TextCell "A108" "7544"
    Set-TextCell "B108" "10/20/2025"
    Set-TextCell "C108" "MONTES DE OCA, MANUEL AV. 1079"
    Set-TextCell "D108" "4"
    Set-TextCell "E108" "810398927"
    Set-TextCell "F108" "PEBCOM"
    Set-TextCell "G108" "Pendiente"
    Set-TextCell "H108" "Picada"
    $ws.Range("I108").Value = 1
    Set-TextCell "J108" "Cambio"
    Set-TextCell "K108" "Sin equipos"
    Set-TextCell "L108" "Pasante"
    $ws.Range("M108").Value = -58.374368
    $ws.Range("N108").Value = -34.640512
    Set-TextCell "O108" "San Telmo"
    Set-TextCell "P108" "Capital Sur"
    Set-TextCell "Q108" "CON-A"
    Set-TextCell "R108" "Fuera de Poligono OVL"

    # Row 109
    Set-TextCell "A109" "7550"
    Set-TextCell "B109" "10/20/2025"
    Set-TextCell "C109" "BROWN, ALTE. AV. 1375"
    Set-TextCell "D109" "4"
    Set-TextCell "E109" "810398930"
    Set-TextCell "F109" "PEBCOM"
    Set-TextCell "G109" "Pendiente"
    Set-TextCell "H109" "Picada"
    $ws.Range("I109").Value = 1
    Set-TextCell "J109" "Cambio"
    Set-TextCell "K109" "Sin equipos"
    Set-TextCell "L109" "Pasante"
    $ws.Range("M109").Value = -58.358182
    $ws.Range("N109").Value = -34.636697
    Set-TextCell "O109" "San Telmo"
    Set-TextCell "P109" "Capital Sur"
    Set-TextCell "Q109" "CON-G"
    Set-TextCell "R109" "Fuera de Poligono OVL"

    # Row 110
    Set-TextCell "A110" "7551"
    Set-TextCell "B110" "10/20/2025"
    Set-TextCell "C110" "BROWN, ALTE. AV. 1405"
    Set-TextCell "D110" "4"
    Set-TextCell "E110" "810398933"
    Set-TextCell "F110" "PEBCOM"
    Set-TextCell "G110" "Pendiente"
    Set-TextCell "H110" "Picada"
    $ws.Range("I110").Value = 1
    Set-TextCell "J110" "Cambio"
    Set-TextCell "K110" "Sin equipos"
    Set-TextCell "L110" "Pasante"
    $ws.Range("M110").Value = -58.35791
    $ws.Range("N110").Value = -34.637135
    Set-TextCell "O110" "San Telmo"
    Set-TextCell "P110" "Capital Sur"
    Set-TextCell "Q110" "CON-E"
    Set-TextCell "R110" "Fuera de Poligono OVL"

    # Row 111
    Set-TextCell "A111" "7553"
    Set-TextCell "B111" "10/20/2025"
    Set-TextCell "C111" "CORRIENTES AV. 4515"
    Set-TextCell "D111" "5"
    Set-TextCell "E111" "810398934"
    Set-TextCell "F111" "PEBCOM"
    Set-TextCell "G111" "Pendiente"
    Set-TextCell "H111" "Picada"
    $ws.Range("I111").Value = 1
    Set-TextCell "J111" "Cambio"
    Set-TextCell "K111" "Sin equipos"
    Set-TextCell "L111" "Pasante"
    $ws.Range("M111").Value = -58.428907
    $ws.Range("N111").Value = -34.602301
    Set-TextCell "O111" "Almagro"
    Set-TextCell "P111" "Capital Sur"
    Set-TextCell "Q111" "CLI-N"
    Set-TextCell "R111" "Fuera de Poligono OVL"

